$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Remove the empty "Title 25" placeholder shape ---
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Title 25") {
        $sh.Delete()
    }
}

# --- Reposition the remaining rectangles (same translate for all three) ---
$pointsPerEmu = 1 / 12700.0

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Rectangle 12") {
        $sh.Left = 3196633 * $pointsPerEmu
        $sh.Top = 2564904 * $pointsPerEmu
    } elseif ($sh.Name -eq "Rectangle 29") {
        $sh.Left = 3196633 * $pointsPerEmu
        $sh.Top = 2881372 * $pointsPerEmu
    } elseif ($sh.Name -eq "Rectangle 4") {
        $sh.Left = 3196633 * $pointsPerEmu
        $sh.Top = 3250704 * $pointsPerEmu
    }
}
